# Applies the commit's change:
#  - rename the existing sheet ID_7e0f17b -> ID_867f88d
#  - update its row 2 data (new product: "Mens Casual Slim Fit")
#  - add a new sheet ID_03f327c with the same header row + a second
#    product row ("Mens Cotton Jacket")

$wb = $excel.ActiveWorkbook

# --- Rename existing (first) sheet and fix its data ---------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ID_867f88d"

$ws1.Range("A2").Value = 15.99
$ws1.Range("B2").Value = "18/02/2025"
$ws1.Range("C2").Value = "Mens Casual Slim Fit"
$ws1.Range("D2").Value = "The color could be slightly different between on the screen and in practice. / Please note that body builds vary by person, therefore, detailed size information should be reviewed below on the product description."
$ws1.Range("E2").Value = "ID_867f88d"

# --- Add the new second sheet (placed right after the first sheet) ------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ID_03f327c"

$ws2.Range("A1").Value = "price"
$ws2.Range("B1").Value = "date"
$ws2.Range("C1").Value = "name"
$ws2.Range("D1").Value = "description"
$ws2.Range("E1").Value = "id"

$ws2.Range("A2").Value = 55.99
$ws2.Range("B2").Value = "18/02/2025"
$ws2.Range("C2").Value = "Mens Cotton Jacket"
$ws2.Range("D2").Value = "great outerwear jackets for Spring/Autumn/Winter, suitable for many occasions, such as working, hiking, camping, mountain/rock climbing, cycling, traveling or other outdoors. Good gift choice for you or your family member. A warm hearted love to Father, husband or son in this thanksgiving or Christmas Day."
$ws2.Range("E2").Value = "ID_03f327c"

# Header style (bold, centered, thin border) should match sheet1's row 1 —
# copy the formatting across so the new sheet's header looks the same.
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats
